$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Defect logs ")
$ws2 = $wb.Worksheets.Item("Types of defects ")

# --- Fill in the three new defect rows on the "Defect logs " sheet ---

# Row 4 - Windowlifter.c / Enviroment / Design / Code Review
$ws1.Range("B4").Value = "Windowlifter.c"
$ws1.Range("C4").Value = 42188
$ws1.Range("C4").NumberFormat = "d-mmm-yy"
$ws1.Range("D4").Value = 0
$ws1.Range("E4").Value = "Enviroment"
$ws1.Range("F4").Value = "Design"
$ws1.Range("G4").Value = "Code Review"
$ws1.Range("H4").Value = 5
$ws1.Range("I4").Value = "The uC is not respecting the order of the functions and is executing in aleatory order."
$ws1.Range("I4").WrapText = $true

# Row 5 - Windowlifter.h / Data / Design / Code Review
$ws1.Range("B5").Value = "Windowlifter.h"
$ws1.Range("C5").Value = 42188
$ws1.Range("C5").NumberFormat = "d-mmm-yy"
$ws1.Range("D5").Value = 1
$ws1.Range("E5").Value = "Data"
$ws1.Range("F5").Value = "Design"
$ws1.Range("G5").Value = "Code Review"
$ws1.Range("H5").Value = 1
$ws1.Range("I5").Value = "Duplicate declaration, variable and function."
$ws1.Range("I5").WrapText = $true

# Row 6 - main.c / Documentation / Design / Code Review
$ws1.Range("B6").Value = "main.c"
$ws1.Range("C6").Value = 42188
$ws1.Range("C6").NumberFormat = "d-mmm-yy"
$ws1.Range("D6").Value = 2
$ws1.Range("E6").Value = "Documentation"
$ws1.Range("F6").Value = "Design"
$ws1.Range("G6").Value = "Code Review"
$ws1.Range("H6").Value = 1
$ws1.Range("I6").Value = "Added c templated and function descriptions."
$ws1.Range("I6").WrapText = $true

# Match the row heights produced by wrapping the long descriptions
$ws1.Rows.Item(4).RowHeight = 45
$ws1.Rows.Item(5).RowHeight = 30
$ws1.Rows.Item(6).RowHeight = 30

# Apply the same wrap-text format to the rest of column I (rows 7-32) so the
# whole description column shares one consistent style
$ws1.Range("I7:I32").WrapText = $true

# Best-fit column E now that it holds the longer "Enviroment"/"Documentation" text
$ws1.Range("E1").EntireColumn.AutoFit()

# Switch the page to portrait orientation
$ws1.PageSetup.Orientation = 1

# --- Update sheet selections / active tab ---
# "Types of defects " was the active tab before; reset its selection and
# make "Defect logs " the active sheet with C6 selected instead.
$ws2.Activate()
$ws2.Range("A1").Select()

$ws1.Activate()
$ws1.Range("C6").Select()
